$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, pushing the existing rows 176:271 down to 177:272.
$ws.Rows("176:176").Insert()

# Populate the newly inserted row 176 with its data.
$ws.Range("A176").Value = 3
$ws.Range("B176").Value = "Femacal de La Calera"
$ws.Range("C176").Value = "Coquimbo"
$ws.Range("D176").Value = 44572
$ws.Range("E176").Value = 5
$ws.Range("F176").Value = 100112043
$ws.Range("G176").Value = "Pepino ensalada"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 110
$ws.Range("K176").Value = 10000
$ws.Range("L176").Value = 11000
$ws.Range("M176").Value = 10545
$ws.Range("N176").Value = "$/caja 70 unidades"
$ws.Range("O176").Value = "Limache"
$ws.Range("P176").Value = 151
$ws.Range("Q176").Value = 70
$ws.Range("R176").Value = "Hortaliza"
